$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Weekly Quantity" ---
$ws1 = $wb.Worksheets.Item("Weekly Quantity")

# Delete rows 4 through 9 (old data no longer needed)
$ws1.Rows("4:9").Delete()

# Update row 3 values: date -> 45102.99999999999, quantity -> 200
$ws1.Cells.Item(3, 1).Value = 45102.99999999999
$ws1.Cells.Item(3, 2).Value = 200

# --- Sheet 2: "Monthly Trend" ---
$ws2 = $wb.Worksheets.Item("Monthly Trend")

# Update row 3 quantity: 570 -> 200
$ws2.Cells.Item(3, 2).Value = 200

# Delete row 4 (old data no longer needed)
$ws2.Rows("4:4").Delete()
